$d = $word.ActiveDocument

# Locate the paragraph that ends the bibliography ("Thomson Pioneira (2008).").
# Everything from the paragraph right after it, through the "© 2020 ..."
# footer paragraph, is the boilerplate footer block that was removed from
# the page (the trailing blank paragraph right before the page break stays).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Thomson Pioneira (2008)." + [char]13) {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    # Walk forward: blank paragraph, "Ver no Jupiter..." paragraph,
    # "© 2020 ..." paragraph -- these three are the ones to remove.
    $blank = $anchor.Next()
    $jupiter = $blank.Next()
    $copyright = $jupiter.Next()

    $start = $anchor.Range.End
    $end = $copyright.Range.End

    $d.Range($start, $end).Delete()
}
